# Release Update 2.0 lebih rapih
#
# - Update the "Tanggal Input" timestamp on row 2 (A2)
# - Replace row 3's submission with a new one ("ayaka")
# - Remove row 4 (the "Aditya" submission) entirely
# - Narrow columns B, C and E slightly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: just the input timestamp moves ---------------------------------
$ws.Range("A2").Value = "2025-05-20 16:48:50"

# --- Row 3: brand-new submission data ---------------------------------------
$ws.Range("B3").Value = "ayaka"

# C3/D3/N3/O3/P3 look numeric/date-like to Excel's auto-detection; format as
# Text first so they land as literal strings, same as the rest of the sheet.
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "27713123"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2025-05-02"

$ws.Range("A3").Value = "2025-05-20 16:53:45"
$ws.Range("E3").Value = "ayaka@hooy.co.id"
$ws.Range("F3").Value = "INZ"
$ws.Range("G3").Value = "TY"
$ws.Range("H3").Value = "1 Tahun"
# I3 (Platform Type = GCP) is unchanged
$ws.Range("J3").Value = "Testing"
$ws.Range("K3").Value = "Physical Server"
$ws.Range("L3").Value = "terserah"
$ws.Range("M3").Value = "terserah"

$ws.Range("N3").NumberFormat = "@"
$ws.Range("N3").Value = "2"

$ws.Range("O3").NumberFormat = "@"
$ws.Range("O3").Value = "16"

$ws.Range("P3").NumberFormat = "@"
$ws.Range("P3").Value = "50"

$ws.Range("Q3").Value = "B"
# R3 (OS Platform = Linux) is unchanged
$ws.Range("S3").Value = "ZIB"
$ws.Range("T3").Value = ""
$ws.Range("U3").Value = "ZIBPRTGRTG10"

# --- Row 4: delete it entirely, shrinking the used range to A1:U3 ----------
$ws.Rows.Item(4).Delete()

# --- Column widths: B, C, E get narrower ------------------------------------
# Range.ColumnWidth is expressed in "characters"; the engine quantizes it to
# the nearest 1/6 character before converting to the raw XML width (raw =
# ColumnWidth + 5/6), so we solve for the ColumnWidth that lands closest to
# each target raw width.
$ws.Columns.Item(2).ColumnWidth = 15.966666666666667   # -> raw width 16.8
$ws.Columns.Item(3).ColumnWidth = 11.166666666666666   # -> raw width 12
$ws.Columns.Item(5).ColumnWidth = 20.766666666666667   # -> raw width 21.6
